$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '42.427.07'
$ws.Range('E2').Value = '  +3.04%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.233.81'
$ws.Range('E3').Value = '  +2.45%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range('E4').Value = '  -0.16%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '251.85'
$ws.Range('E5').Value = '  -0.96%  '

# Row 6: 'XRP' -> 'XRP'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.631'
$ws.Range('E6').Value = '  +0.71%  '

# Row 7: 'Solana' -> 'Solana'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '69.11'
$ws.Range('E7').Value = '  +1.88%  '

# Row 8: 'USDC' -> 'USDC'
$ws.Range('E8').Value = '  -0.14%  '

# Row 9: 'Cardano' -> 'Cardano'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.630'
$ws.Range('E9').Value = '  +9.15%  '

# Row 10: 'Avalanche' -> 'Avalanche'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '40.13'
$ws.Range('E10').Value = '  +7.03%  '

# Row 11: 'OKB' -> 'OKB'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '59.35'
$ws.Range('E11').Value = '  +1.16%  '

# Row 12: 'Dogecoin' -> 'Dogecoin'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0949'
$ws.Range('E12').Value = '  +1.80%  '

# Row 13: 'Polkadot' -> 'Polkadot'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.13'
$ws.Range('E13').Value = '  +0.08%  '

# Row 14: 'TRON' -> 'TRON'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.105'
$ws.Range('E14').Value = '  +0.08%  '

# Row 15: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.567.89'
$ws.Range('E15').Value = '  +2.27%  '

# Row 16: 'Polygon' -> 'Polygon'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.883'
$ws.Range('E16').Value = '  +1.60%  '

# Row 17: 'Chainlink' -> 'Chainlink'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.67'
$ws.Range('E17').Value = '  +1.32%  '

# Row 18: 'WrappedEther' -> 'WrappedEther'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.237.80'
$ws.Range('E18').Value = '  +2.18%  '

# Row 19: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '42.304.39'
$ws.Range('E19').Value = '  +2.77%  '

# Row 20: 'ShibaInu' -> 'ShibaInu'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0969'
$ws.Range('E20').Value = '  +1.85%  '

# Row 21: 'Uniswap' -> 'Uniswap'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.20'
$ws.Range('E21').Value = '  +0.41%  '

# Row 22: 'Litecoin' -> 'Litecoin'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '72.76'
$ws.Range('E22').Value = '  +1.25%  '

# Row 23: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '232.02'
$ws.Range('E23').Value = '  -0.03%  '

# Row 24: 'ImmutableX' -> 'ImmutableX'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.08'
$ws.Range('E24').Value = '  +2.19%  '

# Row 25: 'WEMIXToken' -> 'WEMIXToken'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.94'
$ws.Range('E25').Value = '  +0.65%  '

# Row 26: 'Cosmos' -> 'Cosmos'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.45'
$ws.Range('E26').Value = '  -3.05%  '

# Row 27: 'Dai' -> 'Dai'
$ws.Range('E27').Value = '  +0.09%  '

# Row 28: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.42'
$ws.Range('E28').Value = '  -3.96%  '

# Row 29: 'LEO' -> 'LEO'
$ws.Range('E29').Value = '  -1.66%  '

# Row 30: 'Toncoin' -> 'Toncoin'
$ws.Range('E30').Value = '  +1.69%  '

# Row 31: 'Monero' -> 'Monero'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '167.57'
$ws.Range('E31').Value = '  -0.88%  '

# Row 32: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '20.73'
$ws.Range('E32').Value = '  +0.72%  '

# Row 33: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.11'

# Row 34: 'Kaspa' -> 'Kaspa'
$ws.Range('E34').Value = '  +4.53%  '

# Row 35: 'Hedera' -> 'Hedera'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0780'
$ws.Range('E35').Value = '  +5.44%  '

# Row 36: 'Stellar' -> 'Stellar'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.124'
$ws.Range('E36').Value = '  +0.44%  '

# Row 37: 'InjectiveProtocol' -> 'InjectiveProtocol'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '27.88'
$ws.Range('E37').Value = '  +0.58%  '

# Row 38: 'Filecoin' -> 'Filecoin'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.64'
$ws.Range('E38').Value = '  +0.98%  '

# Row 39: 'RenderToken' -> 'RenderToken'
$ws.Range('E39').Value = '  -2.36%  '

# Row 40: 'VeChain' -> 'VeChain'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0320'
$ws.Range('E40').Value = '  +7.16%  '

# Row 41: 'LidoDAOToken' -> 'LidoDAOToken'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.26'
$ws.Range('E41').Value = '  +2.77%  '

# Row 42: 'Celestia' -> 'Celestia'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '12.58'
$ws.Range('E42').Value = '  -1.16%  '

# Row 43: 'THORChain' -> 'THORChain'
$ws.Range('E43').Value = '  +0.65%  '

# Row 44: 'FTXToken' -> 'MultiversX'
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '63.30'
$ws.Range('E44').Value = '  -1.24%  '

# Row 45: 'MultiversX' -> 'FTXToken'
$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '5.06'
$ws.Range('E45').Value = '  +1.36%  '

# Row 46: 'Algorand' -> 'Algorand'
$ws.Range('E46').Value = '  -0.58%  '

# Row 47: 'FraxShare' -> 'FraxShare'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.70'
$ws.Range('E47').Value = '  +1.10%  '

# Row 48: 'Cronos' -> 'Cronos'
$ws.Range('E48').Value = '  +0.53%  '

# Row 49: 'BinanceUSD' -> 'BinanceUSD'
$ws.Range('E49').Value = '  -0.18%  '

# Row 50: 'ARBITRUM' -> 'ARBITRUM'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.17'
$ws.Range('E50').Value = '  +3.36%  '

# Row 51: 'TrustWalletToken' -> 'SynthetixNetwork'
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.35'
$ws.Range('E51').Value = '  +1.42%  '
